$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness (column C) values per the diff:
# Rows 2-16 (Generation 0-14): 7293 -> 7320
$ws.Range("C2:C16").Value = 7320

# Rows 17-64 (Generation 15-62): 7293 -> 7295
$ws.Range("C17:C64").Value = 7295
